$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.01066386480764514
$ws.Range("D2").Value = 0.04191488182279812
$ws.Range("E2").Value = 0.4111038076194404
$ws.Range("F2").Value = 1.532045218073179
$ws.Range("G2").Value = 0.002451759371712802
$ws.Range("I2").Value = 1.175047952375223
$ws.Range("K2").Value = 1.942888141309311

$ws.Range("C3").Value = 0.01077484301384679
$ws.Range("D3").Value = 0.04255721576186389
$ws.Range("E3").Value = 0.3582844269763967
$ws.Range("F3").Value = 1.485308932098945
$ws.Range("G3").Value = 0.002457872350537738
$ws.Range("I3").Value = 1.13991178790657
$ws.Range("K3").Value = 1.740295478582425

$ws.Range("C4").Value = 0.01084515869750113
$ws.Range("D4").Value = 0.04298746569473622
$ws.Range("E4").Value = 0.3260051226716314
$ws.Range("F4").Value = 1.457932004874507
$ws.Range("G4").Value = 0.002461813168438446
$ws.Range("I4").Value = 1.119356790374894
$ws.Range("K4").Value = 1.616615238086467

$ws.Range("C5").Value = 0.01087436326625379
$ws.Range("D5").Value = 0.04317168355319012
$ws.Range("E5").Value = 0.3128858927219653
$ws.Range("F5").Value = 1.447102684710856
$ws.Range("G5").Value = 0.002463466405241862
$ws.Range("I5").Value = 1.111232921607638
$ws.Range("K5").Value = 1.56639005467332

$ws.Range("C6").Value = 0.01087924601510259
$ws.Range("D6").Value = 0.04320280584328628
$ws.Range("E6").Value = 0.3107094616744348
$ws.Range("F6").Value = 1.445324095898926
$ws.Range("G6").Value = 0.002463743787549118
$ws.Range("I6").Value = 1.109899099304855
$ws.Range("K6").Value = 1.558060698696806

$ws.Range("C7").Value = 0.01084555032691448
$ws.Range("D7").Value = 0.04298991430433574
$ws.Range("E7").Value = 0.3258280554744459
$ws.Range("F7").Value = 1.45778463912022
$ws.Range("G7").Value = 0.00246183527270055
$ws.Range("I7").Value = 1.119246211720935
$ws.Range("K7").Value = 1.615937177735134

$ws.Range("C8").Value = 0.01070168072198285
$ws.Range("D8").Value = 0.04212883945146473
$ws.Range("E8").Value = 0.3928580752183137
$ws.Range("F8").Value = 1.515653825109212
$ws.Range("G8").Value = 0.002453828343391752
$ws.Range("I8").Value = 1.16271930529588
$ws.Range("K8").Value = 1.87288402261737

$ws.Range("C9").Value = 0.01043667230329959
$ws.Range("D9").Value = 0.04073055761361388
$ws.Range("E9").Value = 0.5256697616410264
$ws.Range("F9").Value = 1.639826365754715
$ws.Range("G9").Value = 0.002439605078357626
$ws.Range("I9").Value = 1.256226908164152
$ws.Range("K9").Value = 2.382605797238398

$ws.Range("C10").Value = 0.01025225542309549
$ws.Range("D10").Value = 0.03988824544280689
$ws.Range("E10").Value = 0.6243119278589973
$ws.Range("F10").Value = 1.737894194054519
$ws.Range("G10").Value = 0.002430044025509625
$ws.Range("I10").Value = 1.330210189419418
$ws.Range("K10").Value = 2.760984054436847

$ws.Range("C11").Value = 0.0101705723699812
$ws.Range("D11").Value = 0.03954699608346246
$ws.Range("E11").Value = 0.6694700684629566
$ws.Range("F11").Value = 1.78406339497883
$ws.Range("G11").Value = 0.002425884767610336
$ws.Range("I11").Value = 1.365069355474674
$ws.Range("K11").Value = 2.93403624197822

$ws.Range("C12").Value = 0.01013995829786118
$ws.Range("D12").Value = 0.03942395186782477
$ws.Range("E12").Value = 0.6866154689563757
$ws.Range("F12").Value = 1.801776105178732
$ws.Range("G12").Value = 0.002424336895469725
$ws.Range("I12").Value = 1.378447119129078
$ws.Range("K12").Value = 2.999705235214151

$ws.Range("C13").Value = 0.01014653745977689
$ws.Range("D13").Value = 0.03945017437046161
$ws.Range("E13").Value = 0.6829208350818874
$ws.Range("F13").Value = 1.797951064288469
$ws.Range("G13").Value = 0.002424669053087628
$ws.Range("I13").Value = 1.375558022163958
$ws.Range("K13").Value = 2.985556032750083

$ws.Range("C14").Value = 0.01016804736663701
$ws.Range("D14").Value = 0.03953674844006372
$ws.Range("E14").Value = 0.6708797061344285
$ws.Range("F14").Value = 1.785516000584664
$ws.Range("G14").Value = 0.002425756880541845
$ws.Range("I14").Value = 1.366166373142434
$ws.Range("K14").Value = 2.939436078383665

$ws.Range("C15").Value = 0.01018126416374265
$ws.Range("D15").Value = 0.03959058684223038
$ws.Range("E15").Value = 0.6635101455163834
$ws.Range("F15").Value = 1.777929207362945
$ws.Range("G15").Value = 0.002426426735871796
$ws.Range("I15").Value = 1.360436940839094
$ws.Range("K15").Value = 2.911204374967554

$ws.Range("C16").Value = 0.01025763802904578
$ws.Range("D16").Value = 0.03991140333578613
$ws.Range("E16").Value = 0.6213668154149872
$ws.Range("F16").Value = 1.73490873598044
$ws.Range("G16").Value = 0.002430319652225899
$ws.Range("I16").Value = 1.327956646274814
$ws.Range("K16").Value = 2.749693720392486

$ws.Range("C17").Value = 0.01030505642983925
$ws.Range("D17").Value = 0.04011905791393744
$ws.Range("E17").Value = 0.5955891180913966
$ws.Range("F17").Value = 1.708919873294292
$ws.Range("G17").Value = 0.002432756388180953
$ws.Range("I17").Value = 1.308342363422412
$ws.Range("K17").Value = 2.650852642116718

$ws.Range("C18").Value = 0.01033253813941215
$ws.Range("D18").Value = 0.04024243210076861
$ws.Range("E18").Value = 0.5807890821198498
$ws.Range("F18").Value = 1.694118050176911
$ws.Range("G18").Value = 0.00243417583953988
$ws.Range("I18").Value = 1.29717379634171
$ws.Range("K18").Value = 2.594088735147807

$ws.Range("C19").Value = 0.01034187869577252
$ws.Range("D19").Value = 0.04028487606253606
$ws.Range("E19").Value = 0.5757824958761972
$ws.Range("F19").Value = 1.689131366809931
$ws.Range("G19").Value = 0.002434659522976144
$ws.Range("I19").Value = 1.293411599344225
$ws.Range("K19").Value = 2.574884230277632

$ws.Range("C20").Value = 0.01029998714572322
$ws.Range("D20").Value = 0.04009654420317332
$ws.Range("E20").Value = 0.5983304151264406
$ws.Range("F20").Value = 1.711671253127349
$ws.Range("G20").Value = 0.002432495141511198
$ws.Range("I20").Value = 1.310418606530661
$ws.Range("K20").Value = 2.661365419615265

$ws.Range("C21").Value = 0.01016172076572275
$ws.Range("D21").Value = 0.03951115060932864
$ws.Range("E21").Value = 0.6744152244726536
$ws.Range("F21").Value = 1.789162209588966
$ws.Range("G21").Value = 0.002425436623924888
$ws.Range("I21").Value = 1.368920080604241
$ws.Range("K21").Value = 2.952978846152689

$ws.Range("C22").Value = 0.01007320687083002
$ws.Range("D22").Value = 0.03916465304562422
$ws.Range("E22").Value = 0.724405201366352
$ws.Range("F22").Value = 1.841146942945642
$ws.Range("G22").Value = 0.002420981628559044
$ws.Range("I22").Value = 1.408189947884679
$ws.Range("K22").Value = 3.144371112397891

$ws.Range("C23").Value = 0.0101202789009065
$ws.Range("D23").Value = 0.0393462323557543
$ws.Range("E23").Value = 0.6976991067164278
$ws.Range("F23").Value = 1.813277289213005
$ws.Range("G23").Value = 0.002423344934950232
$ws.Range("I23").Value = 1.387134687034774
$ws.Range("K23").Value = 3.042146093847748

$ws.Range("C24").Value = 0.01030227828238317
$ws.Range("D24").Value = 0.04010671023369028
$ws.Range("E24").Value = 0.5970910130341167
$ws.Range("F24").Value = 1.710426920141515
$ws.Range("G24").Value = 0.002432613193194075
$ws.Range("I24").Value = 1.309479601116209
$ws.Range("K24").Value = 2.656612402807127

$ws.Range("C25").Value = 0.01050654929635719
$ws.Range("D25").Value = 0.04107693047810201
$ws.Range("E25").Value = 0.4895706724523876
$ws.Range("F25").Value = 1.605056942695484
$ws.Range("G25").Value = 0.002443295870331384
$ws.Range("I25").Value = 1.230021053029915
$ws.Range("K25").Value = 2.244053925553487
